# Applies the "F column" (想去人数 / interested-count) updates described by the diff.
# The workbook has four sheets: 展览, 演出, 本地生活, 全部类型.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 88
$ws1.Range("F6").Value  = 891
$ws1.Range("F7").Value  = 486
$ws1.Range("F8").Value  = 4805
$ws1.Range("F9").Value  = 4805
$ws1.Range("F12").Value = 173
$ws1.Range("F16").Value = 7758
$ws1.Range("F17").Value = 257
$ws1.Range("F20").Value = 549
$ws1.Range("F21").Value = 1433
$ws1.Range("F22").Value = 1433
$ws1.Range("F24").Value = 6298
$ws1.Range("F25").Value = 2268
$ws1.Range("F30").Value = 6229
$ws1.Range("F31").Value = 152
$ws1.Range("F32").Value = 43
$ws1.Range("F36").Value = 6575
$ws1.Range("F49").Value = 2169

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 238
$ws2.Range("F7").Value = 39

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1459

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1459
$ws4.Range("F6").Value  = 238
$ws4.Range("F7").Value  = 88
$ws4.Range("F10").Value = 486
$ws4.Range("F11").Value = 4805
$ws4.Range("F12").Value = 4805
$ws4.Range("F15").Value = 173
$ws4.Range("F17").Value = 7758
$ws4.Range("F18").Value = 257
$ws4.Range("F20").Value = 549
$ws4.Range("F21").Value = 1433
$ws4.Range("F23").Value = 6298
$ws4.Range("F24").Value = 2268
$ws4.Range("F25").Value = 39
$ws4.Range("F32").Value = 6229
$ws4.Range("F33").Value = 152
$ws4.Range("F35").Value = 43
$ws4.Range("F37").Value = 6575
